# Refresh the live crypto snapshot (Price / Volume(1h) columns) to match
# the latest data pull performed by the scheduled GitHub Actions job.
#
# Price-column values are written with a leading apostrophe when they look
# like plain numbers (e.g. "4.80", "1.00", "0.120") so Excel keeps them as
# literal text instead of silently converting them to doubles, which would
# round/trim the trailing zeros and reformat values such as "0.0613".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.750.45"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.100.24"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'227.14"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'61.83"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("D13").Value = "2.410.92"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'21.95"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "2.105.24"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "38.726.88"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'71.58"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'6.03"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'226.76"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("D26").Value = "'9.64"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "'170.26"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'1.42"
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").Value = "'19.34"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  +8.32%  "
$ws.Range("D32").Value = "'0.120"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'4.56"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'7.12"
$ws.Range("E35").Value = "  +11.27%  "
$ws.Range("D36").Value = "'0.0613"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "'3.48"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'17.97"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("D42").Value = "'101.56"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "1.524.97"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  +7.59%  "
$ws.Range("D45").Value = "'2.81"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").Value = "'7.75"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("D47").Value = "'0.0909"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("E48").Value = "  +4.52%  "
$ws.Range("D49").Value = "'4.18"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "2.298.07"
$ws.Range("E51").Value = "  -0.07%  "
